$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 has the same style pattern the new row needs (4,5,4,4,4 -> A,C,D style4; B style5; E style4)
$ws.Range("A5:E5").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)

# Fill in the new BOM row (row 27) for the Raspberry Pi Pico
$ws.Range("A27").Value = "Raspberry Pi Pico"
$ws.Range("B27").Value = "PICO"
$ws.Range("C27").Value = "RPI_PICO"
$ws.Range("D27").Value = "C7203002"
$ws.Range("E27").Value = ""

# Match the row height used throughout the rest of the BOM table
$ws.Rows.Item(27).RowHeight = 14.7
